# Actualización semanal de precios de Fruta/Hortaliza (Frambuesa)
# Vega Monumental Concepción.
#
# Esta edición reordena los pares de filas (Primera/Segunda calidad)
# que conforman el reporte semanal, de manera que cada fila quede
# asociada a la fecha/semana que le corresponde junto con sus
# respectivos volúmenes y precios (mínimo, máximo, promedio ponderado
# y $/Kg). Solo cambian las columnas D, L, M, N, O, P y S; el resto de
# columnas (Mercado, Región, Producto, Unidad de comercialización,
# Origen, Kg/unidad, etc.) permanecen iguales en todas las filas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fila 2 toma los datos de la fila original 4
$ws.Range("D2").Value = 44617
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 6000
$ws.Range("O2").Value = 7000
$ws.Range("P2").Value = 6500
$ws.Range("S2").Value = 3250

# Fila 3 toma los datos de la fila original 7
$ws.Range("D3").Value = 44574
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 7000
$ws.Range("O3").Value = 8000
$ws.Range("P3").Value = 7500
$ws.Range("S3").Value = 3750

# Fila 4 toma los datos de la fila original 8
$ws.Range("D4").Value = 44574
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 6000
$ws.Range("O4").Value = 6000
$ws.Range("P4").Value = 6000
$ws.Range("S4").Value = 3000

# Fila 5 toma los datos de la fila original 13
$ws.Range("D5").Value = 44532
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 10000
$ws.Range("O5").Value = 10000
$ws.Range("P5").Value = 10000
$ws.Range("S5").Value = 5000

# Fila 6 toma los datos de la fila original 14
$ws.Range("D6").Value = 44532
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 8000
$ws.Range("O6").Value = 8000
$ws.Range("P6").Value = 8000
$ws.Range("S6").Value = 4000

# Fila 7 toma los datos de la fila original 15
$ws.Range("D7").Value = 44559
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 200
$ws.Range("N7").Value = 6000
$ws.Range("O7").Value = 7000
$ws.Range("P7").Value = 6500
$ws.Range("S7").Value = 3250

# Fila 8 toma los datos de la fila original 16
$ws.Range("D8").Value = 44559
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 5000
$ws.Range("O8").Value = 5000
$ws.Range("P8").Value = 5000
$ws.Range("S8").Value = 2500

# Fila 9 toma los datos de la fila original 9
$ws.Range("D9").Value = 44216
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 200
$ws.Range("N9").Value = 3500
$ws.Range("O9").Value = 4000
$ws.Range("P9").Value = 3750
$ws.Range("S9").Value = 1875

# Fila 10 toma los datos de la fila original 10
$ws.Range("D10").Value = 44216
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 100
$ws.Range("N10").Value = 3000
$ws.Range("O10").Value = 3000
$ws.Range("P10").Value = 3000
$ws.Range("S10").Value = 1500

# Fila 11 toma los datos de la fila original 2
$ws.Range("D11").Value = 44195
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 200
$ws.Range("N11").Value = 3000
$ws.Range("O11").Value = 3500
$ws.Range("P11").Value = 3250
$ws.Range("S11").Value = 1625

# Fila 12 toma los datos de la fila original 3
$ws.Range("D12").Value = 44195
$ws.Range("L12").Value = "Segunda"
$ws.Range("M12").Value = 100
$ws.Range("N12").Value = 2500
$ws.Range("O12").Value = 2500
$ws.Range("P12").Value = 2500
$ws.Range("S12").Value = 1250

# Fila 13 toma los datos de la fila original 11
$ws.Range("D13").Value = 44602
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 200
$ws.Range("N13").Value = 6000
$ws.Range("O13").Value = 7000
$ws.Range("P13").Value = 6500
$ws.Range("S13").Value = 3250

# Fila 14 toma los datos de la fila original 12
$ws.Range("D14").Value = 44602
$ws.Range("L14").Value = "Segunda"
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = 5000
$ws.Range("O14").Value = 5000
$ws.Range("P14").Value = 5000
$ws.Range("S14").Value = 2500

# Fila 15 toma los datos de la fila original 5
$ws.Range("D15").Value = 44609
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 6500
$ws.Range("O15").Value = 7000
$ws.Range("P15").Value = 6750
$ws.Range("S15").Value = 3375

# Fila 16 toma los datos de la fila original 6
$ws.Range("D16").Value = 44609
$ws.Range("L16").Value = "Segunda"
$ws.Range("M16").Value = 50
$ws.Range("N16").Value = 6000
$ws.Range("O16").Value = 6000
$ws.Range("P16").Value = 6000
$ws.Range("S16").Value = 3000

